$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "27.498.26"
$ws.Range("E2").Value = "  -0.16%  "

# Row 3
$ws.Range("D3").Value = "1.829.39"
$ws.Range("E3").Value = "  -1.68%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.79%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "331.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.60%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.003"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.86%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4586"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.71%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3823"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.60%  "

# Row 9
$ws.Range("E9").Value = "  +1.04%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07885"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.92%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9735"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.51%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.10"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.36%  "

# Row 13
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.882"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.58%  "

# Row 14
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.817.96"
$ws.Range("E14").Value = "  -2.07%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.050"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.30%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.005"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.83%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "89.16"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.60%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06627"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.27%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.00001027"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.39%  "

# Row 20
$ws.Range("E20").Value = "  +1.81%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.003"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.97%  "

# Row 22
$ws.Range("D22").Value = "27.474.38"
$ws.Range("E22").Value = "  -0.20%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.335"
$ws.Range("D23").Style = "Normal"

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.82"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.01%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.301"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.26%  "

# Row 26
$ws.Range("D26").Value = "2.037.37"
$ws.Range("E26").Value = "  -2.04%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "156.21"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.68%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.42"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.48%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.067"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.95%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.280"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.99%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "118.24"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.48%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9498"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.07%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09332"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.37%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.577"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.06%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.248"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.71%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.332"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.53%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05924"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.48%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02191"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.99%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.163"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.66%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.027"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.51%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5767"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.33%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1835"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.15%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.05"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.49%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.264"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.18%  "

# Row 45
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.01"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.46%  "

# Row 46
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5459"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.61%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.873"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.88%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.06621"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.95%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "110.55"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.57%  "

# Row 50
$ws.Range("B50").Value = "PaxosStandard"
$ws.Range("C50").Value = "https://coinranking.com/coin/B8xT718SbVhhh+paxosstandard-pax"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.003"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.97%  "

# Row 51
$ws.Range("B51").Value = "EOS"
$ws.Range("C51").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.042"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.47%  "
